{"js": "// Load all paragraphs in the document body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 1) Collapse the three long \"CORE COMPETENCIES\" bullet paragraphs into\n//    a single summary line.\n// ---------------------------------------------------------------------\nconst researchIdx = paragraphs.items.findIndex(p =>\n  p.text.indexOf(\"Research and Analytics: Survey Methodology:\") === 0\n);\nconst programmingIdx = paragraphs.items.findIndex(p =>\n  p.text.indexOf(\"Programming and Development: Python:\") === 0\n);\nconst infrastructureIdx = paragraphs.items.findIndex(p =>\n  p.text.indexOf(\"Data Infrastructure: Cloud Platforms:\") === 0\n);\n\nif (researchIdx === -1 || programmingIdx === -1 || infrastructureIdx === -1) {\n  throw new Error(\"Could not locate the CORE COMPETENCIES paragraphs to collapse.\");\n}\n\n// Replace the text of the first of the three paragraphs, then delete the\n// other two so only a single paragraph remains.\nparagraphs.items[researchIdx].insertText(\n  \"Research and Analytics \u2022 Programming and Development \u2022 Data Infrastructure\",\n  Word.InsertLocation.replace\n);\nparagraphs.items[programmingIdx].delete();\nparagraphs.items[infrastructureIdx].delete();\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) Insert a new \"TECHNICAL SKILLS\" section right after the\n//    \"Built comprehensive survey operations platform...\" bullet, and\n//    before the closing \"For a more detailed...\" paragraph.\n// ---------------------------------------------------------------------\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"text\");\nawait context.sync();\n\nconst anchorIdx = paragraphs2.items.findIndex(p =>\n  p.text.indexOf(\"Built comprehensive survey operations platform\") !== -1\n);\n\nif (anchorIdx === -1) {\n  throw new Error(\"Could not locate the anchor paragraph for the new TECHNICAL SKILLS section.\");\n}\n\nlet anchor = paragraphs2.items[anchorIdx];\n\n// Insert all four new paragraphs off the same anchor (each \"after\" insert\n// lands immediately below the anchor), in reverse order, so the final\n// reading order is: TECHNICAL SKILLS, RESEARCH..., PROGRAMMING...,\n// DATA INFRASTRUCTURE... Anchoring every insert on the original (Normal\n// style, no explicit pPr) paragraph means the three body paragraphs\n// naturally inherit the Normal style without stamping an explicit\n// <w:pStyle>, matching the source formatting; only the heading paragraph\n// gets an explicit style applied.\nconst infrastructurePara = anchor.insertParagraph(\n  \"DATA INFRASTRUCTURE Cloud Platforms; Big Data; Databases; Geospatial\",\n  Word.InsertLocation.after\n);\nconst programmingPara = anchor.insertParagraph(\n  \"PROGRAMMING AND DEVELOPMENT Python; JVM Languages; Web Technologies; Database Languages\",\n  Word.InsertLocation.after\n);\nconst researchPara = anchor.insertParagraph(\n  \"RESEARCH AND ANALYTICS Survey Methodology; Statistical Analysis; Geospatial Analysis; Data Visualization\",\n  Word.InsertLocation.after\n);\nconst headingPara = anchor.insertParagraph(\"TECHNICAL SKILLS\", Word.InsertLocation.after);\nheadingPara.style = \"Heading 2\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndex($doc, $needle) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        if ($doc.Paragraphs.Item($i).Range.Text -like \"*$needle*\") {\n            return $i\n        }\n    }\n    return -1\n}\n\n# ---------------------------------------------------------------------\n# 1) Collapse the three long \"CORE COMPETENCIES\" bullet paragraphs into\n#    a single summary line.\n# ---------------------------------------------------------------------\n$researchIdx = Find-ParagraphIndex $d \"Research and Analytics: Survey Methodology:\"\n$programmingIdx = Find-ParagraphIndex $d \"Programming and Development: Python:\"\n$infrastructureIdx = Find-ParagraphIndex $d \"Data Infrastructure: Cloud Platforms:\"\n\nif ($researchIdx -eq -1 -or $programmingIdx -eq -1 -or $infrastructureIdx -eq -1) {\n    throw \"Could not locate the CORE COMPETENCIES paragraphs to collapse.\"\n}\n\n$bullet = [char]0x2022\n$d.Paragraphs.Item($researchIdx).Range.Text = \"Research and Analytics $bullet Programming and Development $bullet Data Infrastructure\"\n\n# Delete the other two paragraphs (delete the higher index first so the\n# lower index stays valid).\n$d.Paragraphs.Item($infrastructureIdx).Range.Delete() | Out-Null\n$d.Paragraphs.Item($programmingIdx).Range.Delete() | Out-Null\n\n# ---------------------------------------------------------------------\n# 2) Insert a new \"TECHNICAL SKILLS\" section right after the\n#    \"Built comprehensive survey operations platform...\" bullet, and\n#    before the closing \"For a more detailed...\" paragraph.\n# ---------------------------------------------------------------------\n$anchorIdx = Find-ParagraphIndex $d \"Built comprehensive survey operations platform\"\n\nif ($anchorIdx -eq -1) {\n    throw \"Could not locate the anchor paragraph for the new TECHNICAL SKILLS section.\"\n}\n\n$anchor = $d.Paragraphs.Item($anchorIdx)\n\n# Insert all four new paragraphs off the same anchor (each InsertParagraphAfter\n# lands immediately below the anchor), in reverse order, so the final\n# reading order is: TECHNICAL SKILLS, RESEARCH..., PROGRAMMING...,\n# DATA INFRASTRUCTURE... Anchoring every insert on the original (Normal\n# style) paragraph means the three body paragraphs naturally inherit the\n# Normal style; only the heading paragraph gets an explicit style applied.\n$anchor.Range.InsertParagraphAfter()\n$d.Paragraphs.Item($anchorIdx + 1).Range.Text = \"DATA INFRASTRUCTURE Cloud Platforms; Big Data; Databases; Geospatial\"\n\n$anchor.Range.InsertParagraphAfter()\n$d.Paragraphs.Item($anchorIdx + 1).Range.Text = \"PROGRAMMING AND DEVELOPMENT Python; JVM Languages; Web Technologies; Database Languages\"\n\n$anchor.Range.InsertParagraphAfter()\n$d.Paragraphs.Item($anchorIdx + 1).Range.Text = \"RESEARCH AND ANALYTICS Survey Methodology; Statistical Analysis; Geospatial Analysis; Data Visualization\"\n\n$anchor.Range.InsertParagraphAfter()\n$d.Paragraphs.Item($anchorIdx + 1).Range.Text = \"TECHNICAL SKILLS\"\n$d.Paragraphs.Item($anchorIdx + 1).Style = \"Heading 2\"\n"}
